$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 626.29999999999995
$ws.Range("K5").Value = 418.2
$ws.Range("K6").Value = 652.29999999999995
